# Natmi following Dr Hou advice
# Rebuild the LR-pairs table for Anpep-Sele: row 2-6 values are recomputed
# (Ligand-expressing cells / Receptor-expressing cells counts, specificity
# scores, edge weights, ...) and 5 additional Target-cluster="M2" rows are
# appended (rows 7-11), pairing each Sending cluster with both ECs and M2
# as Target cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Anpep"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.880375333333333
$ws.Range("H2").Value = 17.641126
$ws.Range("I2").Value = 0.02508458002889762
$ws.Range("J2").Value = 0.02543725130933246
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.112632333333332
$ws.Range("N2").Value = 27.337897
$ws.Range("O2").Value = 0.9981738658344552
$ws.Range("P2").Value = 0.9981738658344552
$ws.Range("Q2").Value = 53.5856983946691
$ws.Range("R2").Value = 482.2712855520219
$ws.Range("S2").Value = 0.02503877222027851
$ws.Range("T2").Value = 0.02539079947563894
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Anpep"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.880375333333333
$ws.Range("H3").Value = 17.641126
$ws.Range("I3").Value = 0.02508458002889762
$ws.Range("J3").Value = 0.02543725130933246
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01667133333333333
$ws.Range("N3").Value = 0.050014
$ws.Range("O3").Value = 0.001826134165544791
$ws.Range("P3").Value = 0.001826134165544791
$ws.Range("Q3").Value = 0.0980336973071111
$ws.Range("R3").Value = 0.882303275764
$ws.Range("S3").Value = 0.0000458078086191125
$ws.Range("T3").Value = 0.00004645183369352097
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Anpep"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 111.7222593333333
$ws.Range("H4").Value = 335.166778
$ws.Range("I4").Value = 0.4765862375093723
$ws.Range("J4").Value = 0.483286699642939
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.112632333333332
$ws.Range("N4").Value = 27.337897
$ws.Range("O4").Value = 0.9981738658344552
$ws.Range("P4").Value = 0.9981738658344552
$ws.Range("Q4").Value = 1018.083872753985
$ws.Range("R4").Value = 9162.754854785866
$ws.Range("S4").Value = 0.475715927098228
$ws.Range("T4").Value = 0.4824041532889677
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Anpep"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 111.7222593333333
$ws.Range("H5").Value = 335.166778
$ws.Range("I5").Value = 0.4765862375093723
$ws.Range("J5").Value = 0.483286699642939
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01667133333333333
$ws.Range("N5").Value = 0.050014
$ws.Range("O5").Value = 0.001826134165544791
$ws.Range("P5").Value = 0.001826134165544791
$ws.Range("Q5").Value = 1.862559026099111
$ws.Range("R5").Value = 16.763031234892
$ws.Range("S5").Value = 0.0008703104111443094
$ws.Range("T5").Value = 0.0008825463539713545
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Anpep"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 67.60291833333333
$ws.Range("H6").Value = 202.808755
$ws.Range("I6").Value = 0.288381390471254
$ws.Range("J6").Value = 0.2924358268665977
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.112632333333332
$ws.Range("N6").Value = 27.337897
$ws.Range("O6").Value = 0.9981738658344552
$ws.Range("P6").Value = 0.9981738658344552
$ws.Range("Q6").Value = 616.0405394320261
$ws.Range("R6").Value = 5544.364854888235
$ws.Range("S6").Value = 0.2878547673614071
$ws.Range("T6").Value = 0.2919017998119273
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Anpep"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 67.60291833333333
$ws.Range("H7").Value = 202.808755
$ws.Range("I7").Value = 0.288381390471254
$ws.Range("J7").Value = 0.2924358268665977
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01667133333333333
$ws.Range("N7").Value = 0.050014
$ws.Range("O7").Value = 0.001826134165544791
$ws.Range("P7").Value = 0.001826134165544791
$ws.Range("Q7").Value = 1.127030785841111
$ws.Range("R7").Value = 10.14327707257
$ws.Range("S7").Value = 0.00052662310984687
$ws.Range("T7").Value = 0.0005340270546704353
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Anpep"
$ws.Range("C8").Value = "Sele"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 39.46603233333334
$ws.Range("H8").Value = 118.398097
$ws.Range("I8").Value = 0.1683547036320518
$ws.Range("J8").Value = 0.1707216505304549
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.112632333333332
$ws.Range("N8").Value = 27.337897
$ws.Range("O8").Value = 0.9981738658344552
$ws.Range("P8").Value = 0.9981738658344552
$ws.Range("Q8").Value = 359.6394423091121
$ws.Range("R8").Value = 3236.754980782009
$ws.Range("S8").Value = 0.1680472653558192
$ws.Range("T8").Value = 0.170409889891623
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Anpep"
$ws.Range("C9").Value = "Sele"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 39.46603233333334
$ws.Range("H9").Value = 118.398097
$ws.Range("I9").Value = 0.1683547036320518
$ws.Range("J9").Value = 0.1707216505304549
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01667133333333333
$ws.Range("N9").Value = 0.050014
$ws.Range("O9").Value = 0.001826134165544791
$ws.Range("P9").Value = 0.001826134165544791
$ws.Range("Q9").Value = 0.6579513803731112
$ws.Range("R9").Value = 5.921562423358001
$ws.Range("S9").Value = 0.0003074382762326576
$ws.Range("T9").Value = 0.0003117606388318617
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Anpep"
$ws.Range("C10").Value = "Sele"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.7503315
$ws.Range("H10").Value = 19.500663
$ws.Range("I10").Value = 0.04159308835842419
$ws.Range("J10").Value = 0.02811857165067587
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.112632333333332
$ws.Range("N10").Value = 27.337897
$ws.Range("O10").Value = 0.9981738658344552
$ws.Range("P10").Value = 0.9981738658344552
$ws.Range("Q10").Value = 88.85118608761849
$ws.Range("R10").Value = 533.1071165257109
$ws.Range("S10").Value = 0.04151713379872234
$ws.Range("T10").Value = 0.02806722336629825
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Anpep"
$ws.Range("C11").Value = "Sele"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9.7503315
$ws.Range("H11").Value = 19.500663
$ws.Range("I11").Value = 0.04159308835842419
$ws.Range("J11").Value = 0.02811857165067587
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01667133333333333
$ws.Range("N11").Value = 0.050014
$ws.Range("O11").Value = 0.001826134165544791
$ws.Range("P11").Value = 0.001826134165544791
$ws.Range("Q11").Value = 0.162551026547
$ws.Range("R11").Value = 0.975306159282
$ws.Range("S11").Value = 0.00007595455970184173
$ws.Range("T11").Value = 0.00005134828437761839
